# Resolvido problema com dependencias.
# Add a header row (A..F) above the existing data row, and move the
# original numeric values (1..6) down into row 2. Both rows end up
# horizontally centered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

$headers = @("A", "B", "C", "D", "E", "F")
$values  = @(1, 2, 3, 4, 5, 6)

for ($col = 1; $col -le 6; $col++) {
    $headerCell = $ws.Cells.Item(1, $col)
    $headerCell.Value = $headers[$col - 1]
    $headerCell.HorizontalAlignment = $xlCenter

    $dataCell = $ws.Cells.Item(2, $col)
    $dataCell.Value = $values[$col - 1]
    $dataCell.HorizontalAlignment = $xlCenter
}

# Leave the newly entered data row selected, matching the post-edit
# selection state (sqref A2:F2).
$ws.Range("A2:F2").Select()
